# reestruturando as funçoes de criaçao de linhas
#
# Applies to the "Bing" worksheet:
#  - Removes the solid blue background fill from the section-title rows
#    (row 1, row 2 and row 8), keeping their Space Grotesk / bold / white
#    font and centered alignment.
#  - Gives the table header + data rows (rows 3-5 and 9-11, plus the blank
#    spacer cells A7/A13) the "Space Grotesk" font, centered horizontally
#    and vertically.
#  - Gives the "Variaçao" rows (row 6 and row 12) a bold font, centered
#    horizontally and vertically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bing")

# --- Section title rows: keep font/bold/color, just drop the blue fill ---
# (only the merged range's anchor cell - A1/A2/A8 - actually holds a cell
# record; touch only that cell so we don't materialize new cells for the
# rest of the merged area.)
$titleRanges = @("A1", "A2", "A8")
foreach ($addr in $titleRanges) {
    $r = $ws.Range($addr)
    $name = $r.Font.Name
    $bold = $r.Font.Bold
    $color = $r.Font.Color
    $hAlign = $r.HorizontalAlignment
    $vAlign = $r.VerticalAlignment

    $r.ClearFormats()

    $r.Font.Name = $name
    $r.Font.Bold = $bold
    $r.Font.Color = $color
    $r.HorizontalAlignment = $hAlign
    $r.VerticalAlignment = $vAlign
}

# --- Table header + data rows: Space Grotesk, centered ---
$tableRanges = @("A3:F5", "A9:F11", "A7", "A13")
foreach ($addr in $tableRanges) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Space Grotesk"
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}

# --- "Variaçao" rows: bold, centered ---
$varRanges = @("A6:F6", "A12:F12")
foreach ($addr in $varRanges) {
    $r = $ws.Range($addr)
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}
